$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'58.859.26"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.51%  "
$c = $ws.Range("D3")
$c.Value = "'2.493.78"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$c.Value = "'532.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.66%  "
$c = $ws.Range("D6")
$c.Value = "'134.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.48%  "
$ws.Range("E7").Value = "  +0.09%  "
$c = $ws.Range("D8")
$c.Value = "'0.565"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.53%  "
$c = $ws.Range("D9")
$c.Value = "'2.517.52"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.29%  "
$c = $ws.Range("D10")
$c.Value = "'0.0992"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.09%  "
$c = $ws.Range("D11")
$c.Value = "'0.153"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.42%  "
$c = $ws.Range("D12")
$c.Value = "'5.23"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "
$c = $ws.Range("D13")
$c.Value = "'0.333"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "
$c = $ws.Range("D14")
$c.Value = "'2.944.99"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "
$c = $ws.Range("D15")
$c.Value = "'58.965.54"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.84%  "
$c = $ws.Range("D16")
$c.Value = "'22.38"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.70%  "
$c = $ws.Range("D17")
$c.Value = "'0.0000136"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.97%  "
$c = $ws.Range("D18")
$c.Value = "'2.525.37"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.83%  "
$c = $ws.Range("D19")
$c.Value = "'10.67"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("E20").Value = "  +4.20%  "
$c = $ws.Range("D21")
$c.Value = "'321.44"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.59%  "
$c = $ws.Range("D22")
$c.Value = "'6.16"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +9.67%  "
$c = $ws.Range("D23")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$c = $ws.Range("D24")
$c.Value = "'65.77"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.30%  "
$ws.Range("E25").Value = "  +2.18%  "
$c = $ws.Range("D26")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$c = $ws.Range("D27")
$c.Value = "'0.160"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.20%  "
$c = $ws.Range("D28")
$c.Value = "'7.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.09%  "
$c = $ws.Range("D29")
$c.Value = "'0.0₃0764"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +7.10%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D30")
$c.Value = "'1.74"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.30%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D31")
$c.Value = "'170.99"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.55%  "
$c = $ws.Range("D32")
$c.Value = "'1.20"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.35%  "
$c = $ws.Range("D33")
$c.Value = "'6.31"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "
$c = $ws.Range("D34")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$c = $ws.Range("D35")
$c.Value = "'0.994"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "
$c = $ws.Range("D36")
$c.Value = "'18.19"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("E39").Value = "  +5.46%  "
$ws.Range("E40").Value = "  +1.41%  "
$c = $ws.Range("D41")
$c.Value = "'0.786"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D42")
$c.Value = "'278.44"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D43")
$c.Value = "'3.49"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D44")
$c.Value = "'132.05"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +10.25%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D45")
$c.Value = "'5.07"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.75%  "
$c = $ws.Range("D46")
$c.Value = "'0.594"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.14%  "
$c = $ws.Range("D47")
$c.Value = "'0.0933"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.26%  "
$c = $ws.Range("D48")
$c.Value = "'0.0511"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +6.30%  "
$c = $ws.Range("D49")
$c.Value = "'0.0219"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.15%  "
$c = $ws.Range("D50")
$c.Value = "'17.08"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.99%  "
$c = $ws.Range("D51")
$c.Value = "'1.758.71"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.70%  "
